$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.421395778656006
$ws.Range("B1").Value = 2.700477361679077
$ws.Range("C1").Value = 1.871700644493103
$ws.Range("D1").Value = 1.675817489624023
$ws.Range("E1").Value = 1.62127685546875
